$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (merged A2:M2) currently holds the descriptive column headers
# (Make, Model, Serial Number, ...). Capture their text before touching
# anything.
$headers = @()
for ($col = 1; $col -le 13; $col++) {
    $headers += $ws.Cells.Item(2, $col).Value2
}

# Move the descriptive headers down into row 3, which is the Excel table's
# header row - it currently holds the generic placeholder labels
# ("Column1".."Column13"). Writing directly into a ListObject header cell
# makes Excel rename the underlying table column to match the new text, but
# the target file keeps the table's column names (Column1..Column13)
# unchanged, so stage each value in a scratch cell outside the table and
# bring it in via Copy/PasteSpecial (values only), which does not trigger
# that rename.
for ($col = 1; $col -le 13; $col++) {
    $staging = $ws.Cells.Item(500, $col)
    $staging.Value2 = $headers[$col - 1]
    $staging.Copy()
    $ws.Cells.Item(3, $col).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
    $staging.ClearContents()
}

# A2:M2 is a merged range, so individual non-anchor cells (B2:M2) can't be
# written to while merged. Unmerge temporarily to edit the row.
$ws.Range("A2:M2").UnMerge()

# Turn row 2 into a single title cell with the customer/site name.
$ws.Range("B2:M2").ClearContents()
$ws.Range("A2").Value2 = "Rich Products Corporation - HQ @ Buffalo, NY"

# Restore the merge across A2:M2 for the title row.
$ws.Range("A2:M2").Merge()
